$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift rows down: old row4 (headers) -> row5, old row5 (data) -> row6 ---
$ws.Rows("4").Insert()

# ================== ROW 4 : new top-level group headers ==================
$ws.Range("F4").Value2 = "callingCodes"

$ws.Range("S4").Value2 = "currencies"
$ws.Range("S4:U4").Merge()
$ws.Range("S4:U4").Interior.Color = 3329434

$ws.Range("V4").Value2 = "languages"
$ws.Range("V4:Y4").Merge()
$ws.Range("V4:Y4").Interior.Color = 3932402

# ================== ROW 5 : column headers ==================
# callingCodes header moved up to row 4 - clear its old spot
$ws.Range("F5").Value2 = ""

# currencies used to be one JSON blob column (S) - now three scalar columns
$ws.Range("S5").Value2 = "code"
$ws.Range("T5").Value2 = "Name"
$ws.Range("U5").Value2 = "Symbol"

# languages used to be one JSON-ish column (was T, now V) - now four scalar columns
$ws.Range("V5").Value2 = "ISO369_1"
$ws.Range("W5").Value2 = "ISO369_2"
$ws.Range("X5").Value2 = "name"
$ws.Range("Y5").Value2 = "native name"

# the remaining headers slide right to make room, and gain new cells
$ws.Range("Z5").Value2 = "gini"
$ws.Range("AA5").Value2 = "translations"
$ws.Range("AB5").Value2 = "regionalBlocs"
$ws.Range("AC5").Value2 = "Timeone"

# ================== ROW 6 : data ==================
# currency data, split out of the old JSON blob
$ws.Range("S6").Value2 = "MZN"
$ws.Range("T6").Value2 = "Mozambican metical"
$ws.Range("U6").Value2 = "MT"

# language data, split out of the old blob
$ws.Range("V6").Value2 = "pt"
$ws.Range("W6").Value2 = "por"
$ws.Range("X6").Value2 = "Portuguese"
$ws.Range("Y6").Value2 = "Português"

# gini slides right into its own new cell; translations/regionalBlocs data is dropped
$ws.Range("Z6").Value2 = 54
$ws.Range("AA6").Value2 = ""
$ws.Range("AB6").Value2 = ""

# Timeone now holds a real UTC-offset string instead of "System.String[]"
$ws.Range("AC6").Value2 = "[UTC+02:00] "

Write-Output "done"
